$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update perk values (database refresh) ---
# Row 18 (tesla_shield / Uncommon): 0.3-DURATION -> 25-BASE_STAT
$ws.Range("B18").Value = "25-BASE_STAT"
# Row 19 (tesla_shield / Uncommon): -BASE_STAT -> 0.3-DURATION
$ws.Range("B19").Value = "0.3-DURATION"
# Row 15 (nitro / Uncommon): 5%-COOL_DOWN -> 4%-COOL_DOWN
$ws.Range("B15").Value = "4%-COOL_DOWN"

# --- Align formatting of rows 14, 15 and 19 with the rest of the
#     "Uncommon" block (rows 16-18), which use the light fill style ---
$ws.Range("A16").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A16").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A16").Copy()
$ws.Range("A19").PasteSpecial(-4122)

$ws.Range("B16").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B16").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B16").Copy()
$ws.Range("B19").PasteSpecial(-4122)

# --- Update the active cell selection on the sheet ---
$ws.Range("G17").Select() | Out-Null
